# Apply the TestData.xlsx edit: update ListCategoryPage sheet (sheet4)
# Adds a second test-data column (category status) next to the existing
# category column, and fixes the "catgeory expected" typo along the way:
#   A1: "category expected"        A2: "Food"
#   B1: "categoryStatusExpected"   B2: "Active"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ListCategoryPage")

# Set the new column B content first (category status expected/value)
$ws.Range("B1").Value = "categoryStatusExpected"
$ws.Range("B2").Value = "Active"

# Set the new column A content (category expected/value)
$ws.Range("A1").Value = "category expected"
$ws.Range("A2").Value = "Food"

# Set new column width for column B to match the diff (raw OOXML width="22")
# The engine adds ~0.8333 (5/6) character padding on top of ColumnWidth when
# serializing, so back that off here to land exactly on width="22".
$ws.Columns.Item(2).ColumnWidth = 22 - (5/6)

# Update the selected cell on this sheet to C3 (as captured in the saved view state)
$ws.Range("C3").Select()
